# Update countries & provincias Spain
# Refreshes the COVID-19 case numbers for the countries whose stats changed
# and re-sorts the country table by total cases (desc) / new cases (desc),
# matching the published data refresh at 06:17 (previously 05:00).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full data range for the country table (row 4 .. row 219), columns A-H.
$dataRange = $ws.Range("A4:H219")

# New figures: Country -> Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes
$updates = @{
    "India"                  = @(2975701, 2333, 2222577, 697196, 0, 0, 55928)
    "Peru"                   = @(576067,   0,   384908,  163914, 0, 0, 27245)
    "Kazajistan"              = @(104313, 242,   89712,   13186, 0, 0,  1415)
    "Kirguistan"              = @(42703,    0,   35831,    5817, 0, 0,  1055)
    "Jamaica"                 = @(1346,    56,     788,     542, 0, 1,    16)
    "Islas Turcas y Caicos"   = @(347,     13,     102,     243, 0, 0,     2)
    "Mongolia"                = @(298,      0,     288,      10, 0, 0,     0)
    "Butan"                   = @(154,      1,     110,      44, 0, 0,     0)
}

foreach ($country in $updates.Keys) {
    $found = $ws.Range("A4:A219").Find($country)
    if ($found -eq $null) {
        Write-Host "WARNING: country not found: $country"
        continue
    }
    $r = $found.Row
    $vals = $updates[$country]
    $ws.Cells.Item($r, 2).Value2 = $vals[0]
    $ws.Cells.Item($r, 3).Value2 = $vals[1]
    $ws.Cells.Item($r, 4).Value2 = $vals[2]
    $ws.Cells.Item($r, 5).Value2 = $vals[3]
    $ws.Cells.Item($r, 6).Value2 = $vals[4]
    $ws.Cells.Item($r, 7).Value2 = $vals[5]
    $ws.Cells.Item($r, 8).Value2 = $vals[6]
}

# Re-sort the whole table by Casos totales desc, tie-break Nuevos casos desc,
# which is how the source ranking is produced.
$key1 = $ws.Range("B4:B219")
$key2 = $ws.Range("C4:C219")
$dataRange.Sort($key1, 2, $key2, [Type]::Missing, 2, [Type]::Missing, 1, 1)

# Update the "last refreshed" timestamp shown at the top of the sheet.
$ws.Range("A1").Value2 = "Datos actualizados a 22 de Agosto de 2020 a las 06:17"
